$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells holding numeric-looking text values must be forced to Text format
# before assignment, otherwise Excel auto-converts the string to a number
# (losing formatting like trailing zeros) and changes the cell type.
function Set-TextValue($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
}

Set-TextValue $ws "D2" "242.71"
Set-TextValue $ws "D3" "22.97"
Set-TextValue $ws "D4" "5.396"
Set-TextValue $ws "D6" "3.429"
Set-TextValue $ws "D8" "0.8141"
Set-TextValue $ws "D9" "0.9181"
Set-TextValue $ws "D10" "0.1437"
Set-TextValue $ws "D11" "0.07417"
Set-TextValue $ws "D12" "0.03284"
Set-TextValue $ws "D13" "0.03068"
Set-TextValue $ws "D14" "0.09347"
Set-TextValue $ws "D15" "3.863"
Set-TextValue $ws "D16" "0.001582"
Set-TextValue $ws "D17" "0.04712"
Set-TextValue $ws "D18" "0.0005992"
$ws.Range("E18").Value = "17OneONEWorstin24h"
Set-TextValue $ws "D19" "0.005912"
Set-TextValue $ws "D20" "0.001259"
$ws.Range("E20").Value = "19BitKanKANBestin24h"
Set-TextValue $ws "D21" "0.004791"
Set-TextValue $ws "D22" "0.00007998"
Set-TextValue $ws "D23" "3.574"
$ws.Range("E27").Value = "26UpBotsUBXT"
Set-TextValue $ws "D40" "0.03937"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue $ws "D41" "0.006438"
$ws.Range("E41").Value = "40KickTokenKICK"
Set-TextValue $ws "D42" "0.003799"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws "D43" "0.1072"
$ws.Range("E43").Value = "42BKEXTokenBKK"
Set-TextValue $ws "D44" "0.008903"
Set-TextValue $ws "D45" "0.00005179"
Set-TextValue $ws "D47" "0.7002"
Set-TextValue $ws "D48" "0.002143"
Set-TextValue $ws "D49" "0.00002099"
Set-TextValue $ws "D50" "0.0001999"
